$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Update cell M26: new value and style (match style already used by row 27 / s=13)
$ws.Range("M27").Copy()
$ws.Range("M26").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("M26").Value = 1156190

# Update the sheet view: right-to-left, scroll position, and selection
$ws.Application.ActiveWindow.DisplayRightToLeft = $true
$ws.Range("H7").Select()
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Application.ActiveWindow.ScrollColumn = 8
$ws.Range("M25").Select()
